$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# 1) Insert a new row at 21 (pushes old row 20's successors, incl. signature rows 25/26 -> 26/27)
$ws.Rows("21:21").Insert()

# 2) Copy old row 20's formatting (now still at row 20, its "closing border" style) down into new row 21
$ws.Range("B20:J20").Copy($ws.Range("B21:J21"))

# 3) Re-style row 20 as a normal data row (same look as rows 16-19)
$ws.Range("B19:J19").Copy($ws.Range("B20:J20"))

# 4) Update period (column E) values for rows 16-21 to ascending order 2503..2508
$ws.Range("E16").Value = "2503"
$ws.Range("E17").Value = "2504"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2506"
$ws.Range("E20").Value = "2507"
$ws.Range("E21").Value = "2508"

# 5) Update Valor Mora (column F) values for rows 16-21
$ws.Range("F16").Value = 52000
$ws.Range("F17").Value = 56940
$ws.Range("F18").Value = 56940
$ws.Range("F19").Value = 56940
$ws.Range("F20").Value = 56940
$ws.Range("F21").Value = 56940

# 6) Salario Basico (column G) stays the same for all rows - ensure row 21 matches
$ws.Range("G21").Value = 1423500

# 7) Update summary fields
$ws.Range("E11").Value = 336700
$ws.Range("F13").Value = 6

Write-Host "done"
